# Commit: add motor2 information in main page.
#
# The "当前过程" (current-process) row is split into two rows, one per
# motor ("电机1过程" / "电机2过程"), and the "累加角度（脉冲）"
# (accumulated-angle) row is likewise split into "电机1角度" /
# "电机2角度". Every address/row below shifts down to make room, and the
# hex addresses in column A increment by 2 (0x0002) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 is unchanged.
$ws.Cells.Item(2, 1).Value = "0x0000"
$ws.Cells.Item(2, 2).Value = "运行状态"
$ws.Cells.Item(2, 3).Value = "0：停止 1：运行"

# Motor 1 / motor 2 process (replaces the old single "当前过程" row).
$ws.Cells.Item(3, 1).Value = "0x0001"
$ws.Cells.Item(3, 2).Value = "电机1过程"
$ws.Cells.Item(3, 3).Value = "int"

$ws.Cells.Item(4, 1).Value = "0x0002"
$ws.Cells.Item(4, 2).Value = "电机2过程"
$ws.Cells.Item(4, 3).Value = "int"

# Motor 1 / motor 2 angle (replaces the old single "累加角度（脉冲）" row).
$ws.Cells.Item(5, 1).Value = "0x0003"
$ws.Cells.Item(6, 1).Value = "0x0005"
$ws.Cells.Item(7, 1).Value = "0x0007"

$ws.Cells.Item(5, 2).Value = "电机1角度"
$ws.Cells.Item(6, 2).Value = "电机2角度"

$ws.Cells.Item(5, 3).Value = "long int"
$ws.Cells.Item(6, 3).Value = "long int"
$ws.Cells.Item(7, 3).Value = "按键返回"

$ws.Cells.Item(7, 2).Value = "回零按钮"

# The remaining rows keep their names/types but shift down two rows, and
# every address from here on is bumped by 2 (0x0002) to stay in sync.
$data = @"
8	0x0008	运行按钮	按键返回
9	0x000A	停止按钮	按键返回
10	0x000C	前进按钮	按键返回
11	0x000E	后退按钮	按键返回
12	0x0010	跟随按钮	按键返回
13	0x0012	过程一前进按钮\设置按钮	按键返回
14	0x0014	过程一后退按钮	按键返回
15	0x0016	过程二前进按钮	按键返回
16	0x0018	过程二后退按钮	按键返回
17	0x001A	过程三前进按钮	按键返回
18	0x001C	过程三后退按钮	按键返回
19	0x001E	保存按钮	按键返回
20	0x0020	返回按钮	按键返回
21	0x0022	过程一前进	0：非激活 1：激活
22	0x0024	过程一后退	0：非激活 1：激活
23	0x0026	过程二前进	0：非激活 1：激活
24	0x0028	过程二后退	0：非激活 1：激活
25	0x002A	过程三前进	0：非激活 1：激活
26	0x002C	过程三后退	0：非激活 1：激活
27	0x002E	步脉冲数 01	int
28	0x0030	步脉冲数 02	int
29	0x0032	步脉冲数 03	int
30	0x0034	步脉冲数 04	int
31	0x0036	步脉冲数 05	int
32	0x0038	步脉冲数 06	int
33	0x003A	步脉冲数 07	int
34	0x003C	步脉冲数 08	int
35	0x003E	步脉冲数 09	int
36	0x0040	步脉冲数 10	int
37	0x0042	步脉冲数 11	int
38	0x0044	步脉冲数 12	int
39	0x0046	步脉冲数 13	int
40	0x0048	步脉冲数 14	int
41	0x004A	步脉冲数 15	int
42	0x004C	步脉冲数 16	int
43	0x004E	步脉冲数 17	int
44	0x0050	步脉冲数 18	int
45	0x0052	步脉冲数 19	int
46	0x0054	步脉冲数 20	int
47	0x0056	步脉冲数 21	int
48	0x0058	步脉冲数 22	int
49	0x005A	步脉冲数 23	int
50	0x005C	步脉冲数 24	int
51	0x005E	步脉冲数 25	int
52	0x0060	步脉冲数 26	int
53	0x0062	步脉冲数 27	int
54	0x0064	步脉冲数 28	int
55	0x0066	步脉冲数 29	int
56	0x0068	步脉冲数 30	int
57	0x006A	步脉冲数 31	int
58	0x006C	步脉冲数 32	int
59	0x006E	步脉冲数 33	int
60	0x0070	步脉冲数 34	int
61	0x0072	步脉冲数 35	int
62	0x0074	步脉冲数 36	int
63	0x0076	步脉冲数 37	int
64	0x0078	步脉冲数 38	int
65	0x007A	步脉冲数 39	int
66	0x007C	步脉冲数 40	int
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim("`r")
    if ($line -eq "") { continue }
    $parts = $line -split "`t"
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = $parts[1]
    $ws.Cells.Item($r, 2).Value = $parts[2]
    $ws.Cells.Item($r, 3).Value = $parts[3]
}

[void]$ws.Range("D8").Select()
